# Fruta / hortaliza, semanal
# Insert 3 new weekly-report rows at the top of the data block (rows 150-152),
# pushing the existing rows 150-236 down to 153-239.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 blank rows before row 150 (existing data shifts down to 153:239).
$ws.Rows("150:152").Insert()

# --- New row 150: Especial, bandeja 10 kilos ---
$ws.Cells.Item(150, 1).Value = 6
$ws.Cells.Item(150, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(150, 3).Value = "Metropolitana"
$ws.Cells.Item(150, 4).Value = 44806
$ws.Cells.Item(150, 5).Value = 13
$ws.Cells.Item(150, 6).Value = "Fruta"
$ws.Cells.Item(150, 7).Value = 100107
$ws.Cells.Item(150, 8).Value = "Otros"
$ws.Cells.Item(150, 9).Value = 100107002
$ws.Cells.Item(150, 10).Value = "Chirimoya"
$ws.Cells.Item(150, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(150, 12).Value = "Especial"
$ws.Cells.Item(150, 13).Value = 150
$ws.Cells.Item(150, 14).Value = 20000
$ws.Cells.Item(150, 15).Value = 20000
$ws.Cells.Item(150, 16).Value = 20000
$ws.Cells.Item(150, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(150, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(150, 19).Value = 2000
$ws.Cells.Item(150, 20).Value = 10

# --- New row 151: Primera, bandeja 10 kilos ---
$ws.Cells.Item(151, 1).Value = 6
$ws.Cells.Item(151, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(151, 3).Value = "Metropolitana"
$ws.Cells.Item(151, 4).Value = 44806
$ws.Cells.Item(151, 5).Value = 13
$ws.Cells.Item(151, 6).Value = "Fruta"
$ws.Cells.Item(151, 7).Value = 100107
$ws.Cells.Item(151, 8).Value = "Otros"
$ws.Cells.Item(151, 9).Value = 100107002
$ws.Cells.Item(151, 10).Value = "Chirimoya"
$ws.Cells.Item(151, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(151, 12).Value = "Primera"
$ws.Cells.Item(151, 13).Value = 150
$ws.Cells.Item(151, 14).Value = 17000
$ws.Cells.Item(151, 15).Value = 17000
$ws.Cells.Item(151, 16).Value = 17000
$ws.Cells.Item(151, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(151, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(151, 19).Value = 1700
$ws.Cells.Item(151, 20).Value = 10

# --- New row 152: Tercera, bandeja 10 kilos ---
$ws.Cells.Item(152, 1).Value = 6
$ws.Cells.Item(152, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(152, 3).Value = "Metropolitana"
$ws.Cells.Item(152, 4).Value = 44806
$ws.Cells.Item(152, 5).Value = 13
$ws.Cells.Item(152, 6).Value = "Fruta"
$ws.Cells.Item(152, 7).Value = 100107
$ws.Cells.Item(152, 8).Value = "Otros"
$ws.Cells.Item(152, 9).Value = 100107002
$ws.Cells.Item(152, 10).Value = "Chirimoya"
$ws.Cells.Item(152, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(152, 12).Value = "Tercera"
$ws.Cells.Item(152, 13).Value = 150
$ws.Cells.Item(152, 14).Value = 15000
$ws.Cells.Item(152, 15).Value = 15000
$ws.Cells.Item(152, 16).Value = 15000
$ws.Cells.Item(152, 17).Value = "`$/bandeja 10 kilos"
$ws.Cells.Item(152, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(152, 19).Value = 1500
$ws.Cells.Item(152, 20).Value = 10
